# Actualización automática 2025-08-29 13:55:09
# Update August ("agosto") sales figures for client rows 10, 12, 26, 45
# on "VENTAS POR GRUPO", roll the changes into "VENTA MENSUAL" (agosto
# column + yearly total row), and refresh the "CUMPLIMIENTO MENSUAL"
# summary (VENTA / POR CUMPLIR / CUMPLIMIENTO) for the affected product
# groups and the grand total.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO (sheet 1): new sales captured for agosto ---
$wsGrupo.Range("H10").Value = 918.9
$wsGrupo.Range("I10").Value = 387
$wsGrupo.Range("C12").Value = 1555.2
$wsGrupo.Range("H26").Value = 782.1
$wsGrupo.Range("D45").Value = 1418.69

# Row 55 keeps a "N de 53" count of non-zero clients per product group;
# update the counts for the groups that just gained a sale.
$wsGrupo.Range("C55").Value = "4 de 53"
$wsGrupo.Range("H55").Value = "6 de 53"
$wsGrupo.Range("I55").Value = "7 de 53"

# --- VENTA MENSUAL (sheet 2): agosto column (F) totals per client ---
$wsMensual.Range("F10").Value = 1305.9
$wsMensual.Range("F12").Value = 4156.47
$wsMensual.Range("F26").Value = 4090.13
$wsMensual.Range("F45").Value = 3809.67
$wsMensual.Range("F55").Value = 105833.86

# --- CUMPLIMIENTO MENSUAL (sheet 3): VENTA / POR CUMPLIR / CUMPLIMIENTO ---
# Row 2: 240X120 PORCELANATO
$wsCumpl.Range("D2").Value = 5084.47
$wsCumpl.Range("E2").Value = 4885.87304517915
$wsCumpl.Range("F2").Value = 0.509959384241893

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 18258.04
$wsCumpl.Range("E3").Value = 9198.9676
$wsCumpl.Range("F3").Value = 0.6649683121331839

# Row 7: INODOROS
$wsCumpl.Range("D7").Value = 3710.7
$wsCumpl.Range("E7").Value = -1310.7
$wsCumpl.Range("F7").Value = 1.546125

# Row 8: LAVABOS
$wsCumpl.Range("D8").Value = 2404.33
$wsCumpl.Range("E8").Value = -1404.33
$wsCumpl.Range("F8").Value = 2.40433

# Row 19: TOTAL
$wsCumpl.Range("D19").Value = 105833.86
$wsCumpl.Range("E19").Value = 11605.83064517915
$wsCumpl.Range("F19").Value = 0.9011762498570957
